$d = $word.ActiveDocument

# --- Change 1: merge the split "Chapter 5" runs into a single run ---
$paras = $d.Paragraphs
$titlePara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Hyndman and Athanasopoulos*Chapter*5*") {
        $titlePara = $p
    }
}
$tr = $titlePara.Range
$tr.Find.Execute("Notes for Hyndman and Athanasopoulos – Chapter 5", $true, $false, $false, $false, $false, $true, 1, $false, "Notes for Hyndman and Athanasopoulos – Chapter 5", 2) | Out-Null

# --- Change 2: expand the residual-standard-error bullet into the full set of notes ---
$paras = $d.Paragraphs
$targetPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*average error that the model produces*") {
        $targetPara = $p
    }
}
$targetIndex = $targetPara.Index
$trailingEmptyPara = $paras.Item($targetIndex + 1)
$replaceRange = $d.Range($targetPara.Range.Start, $trailingEmptyPara.Range.End)
$xmlPayload = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Another measure of how well the model has fitted the data is the standard deviation of the residuals (the residual standard error).</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> This measure is related to the size of the average error that the model produces.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">The differences between the observed </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/><w:i/></w:rPr><w:t xml:space="preserve">y </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">values and the corresponding fitted values </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/><w:i/></w:rPr><w:t>y-hat</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> are the training-set errors or “residuals”.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> There are a series of plots that serve to perform diagnostics on the residuals:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>ACF plot of residuals to visually detect residual autocorrelation. If autocorrelation is present, then the forecasts are inefficient (there are forecasts that have lower variability).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Histogram of residuals to detect normality.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Residual plot against predictors. We expect the residuals to be randomly scattered without showing any systematic patterns.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> If the scatterplots show a pattern, then the relationship may be </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>nonlinear,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> and the model will need to be modified accordingly.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Residual plot against fitted values: a plot of the residuals against the fitted values should also show no pattern. If a pattern is observed, there may be heteroscedasticity (variance not constant).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">Observations that take extreme values compared to the majority of the data are called </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/><w:b/></w:rPr><w:t>outliers.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> Observations that have a large influence on the estimated coefficients of a regression model are called </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/><w:b/></w:rPr><w:t>influential observations.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">More often than not, time series data are non-stationary: that is, the values of the time series do not fluctuate around a constant mean or with a constant variance. Regressing non-stationary </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:lastRenderedPageBreak/><w:t>time-series can lead to spurious regressions.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> High R-squared and high residual autocorrelation can be signs of spurious regression</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Useful predictors for time-series regression models:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">Trend: it is common for time-series data to be trending. A linear trend can be modeled by simply setting </w:t></w:r><m:oMath><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>x</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>1,t</m:t></m:r></m:sub></m:sSub></m:oMath><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> = t as a predictor.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Dummy variables. Indicator variables that takes the value of 1 (“yes”) or 0 (“no”). The interpretation of the associated coefficient with the dummy variable is that it is a measure of the effect of that category relative to the omitted category.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Intervention variables: it is often necessary to model interventions that may have affected the variable to be forecast.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>When the effect lasts only for one period, we use a “spike” variable. This is a dummy variable that takes the value of one in the period of intervention and zero elsewhere.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>When the level shifts, we use a “step variable”. A step variable takes the value of one after the intervention and zero before.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Trading days in sales data. The number of trading days in each month can be included as a predictor.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Distributed lags, such as the ones that measure the effect of advertising.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Fourier series instead of seasonal dummy variables for long seasonal periods.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Predictive accuracy is a way to determine predictor selection.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Forecasters should not use R-squared to determine whether a model will give good predictions as it will lead to overfitting.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> Therefore, it will always choose the model with most variables.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Adjusted r-squared is a way to select predictors and is equivalent to minimizing the standard error of the</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>regression.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Another method is performing leave-one-out cross validation and compute the mean squared error.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Akaike’s information criterion is an estimator of the relative quality of statistical models for a given set of data.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> The idea is to penalize the fit of the model (SSE) with the number of parameters that need to be estimated.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> For small series, the AIC tends to select too many predictors, so a bias-corrected version is available.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Where possible, all potential regression models should be fitted and the best model should be selected based on one of the measures discussed</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> (this is known as “best subsets” regression).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>When using regression models for time series data, we need to distinguish between the different types of forecasts that can be produced, depending on what is assumed to be known when the forecasts are computed.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">Ex-ante forecasts: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>forecasts made using only the information that is available in advance.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> In order to generate ex-ante forecasts, the model requires forecasts of the predictors.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">Ex-post forecasts: forecasts that are made using later information on the predictors. The model from which ex-post forecasts are produced should not be estimated using </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:lastRenderedPageBreak/><w:t>the data from the forecast period. We assume prior knowledge of the predictor variables (the x variables), but should not assume knowledge of the data that are to be forecast.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>For models that rely on special predictors (seasonal dummies or public holiday indicators, there is no difference between ex-ante and ex-post forecasts</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>, because they rely on predictors known in advance and that are based in calendar variables that repeat themselves.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Scenario based forecasting: in this setting, the forecaster assumes possible scenarios for the predictor variables that are of interest.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>Prediction intervals do not include the uncertainty associated with the future distribution of the predictor variables.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> They assume that the values of the predictors are known in advance.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>The great advantage of regression models is that they can be used to capture important relationships between the forecast variable of interest and the predictor variables.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> However, ex-ante forecasting requires obtaining forecasts of the predictors and that can be challenging.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>An alternative formulation is to use as predictors their lagged values.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> The predictor set is formed by predictor values that are observed </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/><w:i/></w:rPr><w:t xml:space="preserve">h </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">time periods prior to observing </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/><w:i/></w:rPr><w:t>y</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">The simplest way of modelling a nonlinear relationship is to transform the forecast variable y and/or the predictor variable </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>x before estimating the regression model. While this provides a non-linear functional form, the model is still linear in the parameters.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t>There are cases where simply transforming the data will not be adequate and a more general specification is required.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> (We allow f(x) to be a more flexible nonlinear function of x).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve">One of the simplest specifications is to make </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/><w:i/></w:rPr><w:t>f</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> piecewise linear.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> We introduce points where the slope of </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/><w:i/></w:rPr><w:t xml:space="preserve">f </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Roboto Condensed Light" w:hAnsi="Roboto Condensed Light"/></w:rPr><w:t xml:space="preserve"> can change.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$replaceRange.InsertXML($xmlPayload)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
